# Updated cryptos list values (Price / Volume(1h)) per GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.203.21"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "'1.904.22"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'306.67"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.5262"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.3779"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "'0.07252"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").Value = "'21.15"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'0.8999"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'0.08424"
$ws.Range("E12").Value = "  +10.48%  "
$ws.Range("D13").Value = "'1.913.26"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "'94.76"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'5.273"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'27.234.24"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'5.062"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'2.144.23"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'6.437"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'146.84"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "'2.278"
$ws.Range("E26").Value = "  +5.81%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").Value = "'114.88"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'4.921"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "'0.8093"
$ws.Range("E33").Value = "  +6.65%  "
$ws.Range("D34").Value = "'0.05068"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'1.238"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").Value = "'2.952"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("D37").Value = "'3.369"
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("D38").Value = "'2.633"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").Value = "'0.5729"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "'0.01989"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "'1.073"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'6.641"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'8.970"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'117.48"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "'0.1517"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "'0.4844"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "'10.22"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "'37.45"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "'63.66"
$ws.Range("E51").Value = "  +0.05%  "
